$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1171.8723  # was 1224.591
$ws.Range("I19").Value = 992.25  # was 1031.0454
$ws.Range("J19").Value = 1359.3043  # was 1418.1364
$ws.Range("K19").Value = 992.25  # was 1031.0454
$ws.Range("L19").Value = 1359.3043  # was 1418.1364
$ws.Range("M19").Value = -817.25  # was -856.0454
$ws.Range("N19").Value = -1709.3043  # was -1768.1364
# Row 64
$ws.Range("H64").Value = 4874.3  # was 4960.3335
$ws.Range("I64").Value = 5780  # was 6200
$ws.Range("K64").Value = 5780  # was 6200
$ws.Range("M64").Value = -5532  # was -5952
# Row 67
$ws.Range("H67").Value = 4874.3  # was 4960.3335
$ws.Range("I67").Value = 5780  # was 6200
$ws.Range("K67").Value = 5780  # was 6200
$ws.Range("M67").Value = -4922  # was -5342
# Row 74
$ws.Range("H74").Value = 3920452  # was 2549540.8
$ws.Range("I74").Value = 4632561.5  # was 2832395.2
$ws.Range("K74").Value = 4632561.5  # was 2832395.2
$ws.Range("M74").Value = -4631625.5  # was -2831459.2
# Row 76
$ws.Range("H76").Value = 275002500  # was 42310770
$ws.Range("I76").Value = 275002500  # was 45836332
$ws.Range("J76").Value = 0  # was 4000.5
$ws.Range("K76").Value = 275002500  # was 45836332
$ws.Range("L76").Value = 0  # was 4000.5
$ws.Range("M76").Value = -275002185  # was -45836017
$ws.Range("N76").ClearContents()  # was -4630.5
# Row 77
$ws.Range("H77").Value = 3920452  # was 2549540.8
$ws.Range("I77").Value = 4632561.5  # was 2832395.2
$ws.Range("K77").Value = 23162807.5  # was 14161976
$ws.Range("M77").Value = -23158127.5  # was -14157296
# Row 79
$ws.Range("H79").Value = 275002500  # was 42310770
$ws.Range("I79").Value = 275002500  # was 45836332
$ws.Range("J79").Value = 0  # was 4000.5
$ws.Range("K79").Value = 275002500  # was 45836332
$ws.Range("L79").Value = 0  # was 4000.5
$ws.Range("M79").Value = -275001408  # was -45835240
$ws.Range("N79").ClearContents()  # was -6184.5
# Row 135
$ws.Range("H135").Value = 1307.3914  # was 1344.091
$ws.Range("J135").Value = 7750  # was 15000
$ws.Range("L135").Value = 69750  # was 135000
$ws.Range("N135").Value = -74820  # was -140070
# Row 137
$ws.Range("H137").Value = 1053.4706  # was 1061.091
$ws.Range("I137").Value = 928  # was 932.8461
$ws.Range("K137").Value = 2784  # was 2798.5383
$ws.Range("M137").Value = -234  # was -248.5383000000002
# Row 141
$ws.Range("H141").Value = 4577  # was 1864.3611
$ws.Range("I141").Value = 993.3333  # was 655.2258
$ws.Range("J141").Value = 9952.5  # was 9361
$ws.Range("K141").Value = 2979.9999  # was 1965.6774
$ws.Range("L141").Value = 29857.5  # was 28083
$ws.Range("M141").Value = 2200.0001  # was 3214.3226
$ws.Range("N141").Value = -40217.5  # was -38443

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 13222  # was 14222
$ws.Range("J44").Value = 13222  # was 14222
$ws.Range("L44").Value = 13222  # was 14222
$ws.Range("N44").Value = -14198  # was -15198
# Row 74
$ws.Range("H74").Value = 965.2083  # was 933.1667
$ws.Range("I74").Value = 1053.6111  # was 979.8
$ws.Range("K74").Value = 1053.6111  # was 979.8
$ws.Range("M74").Value = -179.6111000000001  # was -105.8
# Row 77
$ws.Range("H77").Value = 965.2083  # was 933.1667
$ws.Range("I77").Value = 1053.6111  # was 979.8
$ws.Range("K77").Value = 5268.0555  # was 4899
$ws.Range("M77").Value = -900.0555000000004  # was -531
# Row 122
$ws.Range("H122").Value = 1001  # was 834.8333
$ws.Range("I122").Value = 1027.25  # was 810.6429000000001
$ws.Range("J122").Value = 966  # was 919.5
$ws.Range("K122").Value = 3081.75  # was 2431.9287
$ws.Range("L122").Value = 2898  # was 2758.5
$ws.Range("M122").Value = -631.75  # was 18.07129999999961
$ws.Range("N122").Value = -7798  # was -7658.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 7830  # was 8052.069
$ws.Range("I105").Value = 6544.9546  # was 7034.95
$ws.Range("J105").Value = 11363.875  # was 10312.333
$ws.Range("K105").Value = 6544.9546  # was 7034.95
$ws.Range("L105").Value = 11363.875  # was 10312.333
$ws.Range("M105").Value = -4797.9546  # was -5287.95
$ws.Range("N105").Value = -14857.875  # was -13806.333
# Row 134
$ws.Range("H134").Value = 82755.32000000001  # was 56153.164
$ws.Range("I134").Value = 3194.4375  # was 2098.4482
$ws.Range("J134").Value = 224196.89  # was 252101.5
$ws.Range("K134").Value = 9583.3125  # was 6295.344599999999
$ws.Range("L134").Value = 672590.67  # was 756304.5
$ws.Range("M134").Value = -7048.3125  # was -3760.344599999999
$ws.Range("N134").Value = -677660.67  # was -761374.5
# Row 137
$ws.Range("H137").Value = 62256  # was 68956
$ws.Range("J137").Value = 62256  # was 68956
$ws.Range("L137").Value = 62256  # was 68956
$ws.Range("N137").Value = -72456  # was -79156

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2659.2195  # was 2797.7632
$ws.Range("I31").Value = 4165.4443  # was 4529.875
$ws.Range("J31").Value = 1480.4348  # was 1538.0454
$ws.Range("K31").Value = 4165.4443  # was 4529.875
$ws.Range("L31").Value = 1480.4348  # was 1538.0454
$ws.Range("M31").Value = -3870.4443  # was -4234.875
$ws.Range("N31").Value = -2070.4348  # was -2128.0454
# Row 34
$ws.Range("H34").Value = 2659.2195  # was 2797.7632
$ws.Range("I34").Value = 4165.4443  # was 4529.875
$ws.Range("J34").Value = 1480.4348  # was 1538.0454
$ws.Range("K34").Value = 4165.4443  # was 4529.875
$ws.Range("L34").Value = 1480.4348  # was 1538.0454
$ws.Range("M34").Value = -3963.4443  # was -4327.875
$ws.Range("N34").Value = -1884.4348  # was -1942.0454
# Row 132
$ws.Range("H132").Value = 2031.6296  # was 2223.0833
$ws.Range("I132").Value = 1633.2727  # was 1469
$ws.Range("J132").Value = 2305.5  # was 3278.8
$ws.Range("K132").Value = 4899.8181  # was 4407
$ws.Range("L132").Value = 6916.5  # was 9836.400000000001
$ws.Range("M132").Value = -2369.8181  # was -1877
$ws.Range("N132").Value = -11976.5  # was -14896.4
# Row 134
$ws.Range("H134").Value = 3132.7551  # was 3432.7556
$ws.Range("I134").Value = 2102.3948  # was 2545.484
$ws.Range("J134").Value = 6692.1816  # was 5397.4287
$ws.Range("K134").Value = 6307.1844  # was 7636.451999999999
$ws.Range("L134").Value = 20076.5448  # was 16192.2861
$ws.Range("M134").Value = -3772.1844  # was -5101.451999999999
$ws.Range("N134").Value = -25146.5448  # was -21262.2861

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1191.6316  # was 1293.1515
$ws.Range("I5").Value = 845.931  # was 910.6087
$ws.Range("J5").Value = 2305.5557  # was 2173
$ws.Range("K5").Value = 2537.793  # was 2731.8261
$ws.Range("L5").Value = 6916.6671  # was 6519
$ws.Range("M5").Value = -2425.793  # was -2619.8261
$ws.Range("N5").Value = -7140.6671  # was -6743
# Row 63
$ws.Range("H63").Value = 114111  # was 4020
$ws.Range("I63").Value = 666.6667  # was 800
$ws.Range("J63").Value = 170833.17  # was 4377.778
$ws.Range("K63").Value = 2000.0001  # was 2400
$ws.Range("L63").Value = 512499.51  # was 13133.334
$ws.Range("M63").Value = -1251.0001  # was -1651
$ws.Range("N63").Value = -513997.51  # was -14631.334
# Row 66
$ws.Range("H66").Value = 114111  # was 4020
$ws.Range("I66").Value = 666.6667  # was 800
$ws.Range("J66").Value = 170833.17  # was 4377.778
$ws.Range("K66").Value = 6000.0003  # was 7200
$ws.Range("L66").Value = 1537498.53  # was 39400.002
$ws.Range("M66").Value = -2256.0003  # was -3456
$ws.Range("N66").Value = -1544986.53  # was -46888.002
# Row 114
$ws.Range("H114").Value = 3195.25  # was 2254.6667
$ws.Range("I114").Value = 2000  # was 116.5
$ws.Range("J114").Value = 3366  # was 6531
$ws.Range("K114").Value = 6000  # was 349.5
$ws.Range("L114").Value = 10098  # was 19593
$ws.Range("M114").Value = -2746  # was 2904.5
$ws.Range("N114").Value = -16606  # was -26101
# Row 115
$ws.Range("H115").Value = 1798  # was 1801
$ws.Range("J115").Value = 1922.5  # was 1926.25
$ws.Range("L115").Value = 5767.5  # was 5778.75
$ws.Range("N115").Value = -8117.5  # was -8128.75
# Row 131
$ws.Range("H131").Value = 5447028  # was 6424557.5
$ws.Range("I131").Value = 38539576  # was 41751184
$ws.Range("J131").Value = 1418.7089  # was 1534.3788
$ws.Range("K131").Value = 115618728  # was 125253552
$ws.Range("L131").Value = 4256.126700000001  # was 4603.136399999999
$ws.Range("M131").Value = -115613688  # was -125248512
$ws.Range("N131").Value = -14336.1267  # was -14683.1364
# Row 135
$ws.Range("H135").Value = 1191.6316  # was 1293.1515
$ws.Range("I135").Value = 845.931  # was 910.6087
$ws.Range("J135").Value = 2305.5557  # was 2173
$ws.Range("K135").Value = 7613.379000000001  # was 8195.478300000001
$ws.Range("L135").Value = 20750.0013  # was 19557
$ws.Range("M135").Value = -5078.379000000001  # was -5660.478300000001
$ws.Range("N135").Value = -25820.0013  # was -24627

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3747.0588  # was 3552.6316
$ws.Range("I80").Value = 3384.6155  # was 3406.25
$ws.Range("J80").Value = 4925  # was 4333.3335
$ws.Range("K80").Value = 3384.6155  # was 3406.25
$ws.Range("L80").Value = 4925  # was 4333.3335
$ws.Range("M80").Value = -2386.6155  # was -2408.25
$ws.Range("N80").Value = -6921  # was -6329.3335
# Row 83
$ws.Range("H83").Value = 3747.0588  # was 3552.6316
$ws.Range("I83").Value = 3384.6155  # was 3406.25
$ws.Range("J83").Value = 4925  # was 4333.3335
$ws.Range("K83").Value = 16923.0775  # was 17031.25
$ws.Range("L83").Value = 24625  # was 21666.6675
$ws.Range("M83").Value = -11931.0775  # was -12039.25
$ws.Range("N83").Value = -34609  # was -31650.6675
# Row 122
$ws.Range("H122").Value = 1317790.1  # was 3291298.5
$ws.Range("I122").Value = 1646687.6  # was 3291298.5
$ws.Range("J122").Value = 2200  # was 0
$ws.Range("K122").Value = 4940062.800000001  # was 9873895.5
$ws.Range("L122").Value = 6600  # was 0
$ws.Range("M122").Value = -4937612.800000001  # was -9871445.5
$ws.Range("N122").Value = -11500  # newly added cell

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 1965.5555  # was 1939.8182
$ws.Range("J136").Value = 3444.0625  # was 3273.8235
$ws.Range("L136").Value = 10332.1875  # was 9821.470499999999
$ws.Range("N136").Value = -15432.1875  # was -14921.4705

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2952.5  # was 1558.5714
$ws.Range("I122").Value = 0  # was 1000
$ws.Range("J122").Value = 2952.5  # was 2303.3333
$ws.Range("K122").Value = 0  # was 3000
$ws.Range("L122").Value = 8857.5  # was 6909.999899999999
$ws.Range("M122").ClearContents()  # was -550
$ws.Range("N122").Value = -13757.5  # was -11809.9999
# Row 132
$ws.Range("H132").Value = 1683.3948  # was 1753.75
$ws.Range("I132").Value = 1433.3334  # was 1480.2609
$ws.Range("J132").Value = 2112.0715  # was 2237.6155
$ws.Range("K132").Value = 4300.0002  # was 4440.7827
$ws.Range("L132").Value = 6336.2145  # was 6712.8465
$ws.Range("M132").Value = -1770.0002  # was -1910.7827
$ws.Range("N132").Value = -11396.2145  # was -11772.8465
